$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new rows for "Alumno C" data block (rows 12-15)
$ws.Range("A12").Value = 11.0
$ws.Range("B12").Value = "El que se duerme pierde"
$ws.Range("C12").Value = "Tom Peter"
$ws.Range("D12").Value = 16.0

$ws.Range("A13").Value = 12.0
$ws.Range("B13").Value = "Sin lugar a duda"
$ws.Range("C13").Value = "Ana Gutierrez"
$ws.Range("D13").Value = 26.0

$ws.Range("A14").Value = 13.0
$ws.Range("B14").Value = "El arte de dormir"
$ws.Range("C14").Value = "Nico"
$ws.Range("D14").Value = 32.0

$ws.Range("A15").Value = 14.0
$ws.Range("B15").Value = "Buscando a Nemo"
$ws.Range("C15").Value = "Humble Po"
$ws.Range("D15").Value = 41.0

# Update existing row 3 (Alumno / book record) - new author and price
$ws.Range("C3").Value = "j. r. tolkien"
$ws.Range("D3").Value = 300.0
